$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Maze20_0.txt)
$ws.Range("B2").Value = 92
$ws.Range("D2").Value = 116
$ws.Range("F2").Value = "40ms"
$ws.Range("G2").Value = "39ms"

# Row 3 (Maze20_1.txt)
$ws.Range("B3").Value = 4
$ws.Range("D3").Value = 4
$ws.Range("F3").Value = "33ms"
$ws.Range("G3").Value = "38ms"

# Row 4 (Maze20_2.txt)
$ws.Range("B4").Value = 237
$ws.Range("D4").Value = 374
$ws.Range("F4").Value = "40ms"
$ws.Range("G4").Value = "40ms"

# Row 5 (Maze20_3.txt)
$ws.Range("G5").Value = "38ms"

# Row 6 (Maze20_4.txt)
$ws.Range("G6").Value = "41ms"

# Row 7 (Maze50_0.txt)
$ws.Range("G7").Value = "47ms"

# Row 8 (Maze50_1.txt)
$ws.Range("C8").Value = 53
$ws.Range("E8").Value = 968
$ws.Range("G8").Value = "43ms"

# Row 9 (Maze50_2.txt)
$ws.Range("C9").Value = 148
$ws.Range("E9").Value = 186
$ws.Range("G9").Value = "43ms"

# Row 10 (Maze50_3.txt)
$ws.Range("C10").Value = 541
$ws.Range("E10").Value = 862
$ws.Range("G10").Value = "41ms"

# Row 11 (Maze50_4.txt)
$ws.Range("C11").Value = 388
$ws.Range("E11").Value = 933
$ws.Range("G11").Value = "43ms"

# Row 12 (Maze100_0.txt)
$ws.Range("C12").Value = 153
$ws.Range("E12").Value = 9925
$ws.Range("G12").Value = "71ms"

# Row 13 (Maze100_1.txt)
$ws.Range("C13").Value = 4369
$ws.Range("E13").Value = 8051
$ws.Range("G13").Value = "79ms"

# Row 14 (Maze100_2.txt)
$ws.Range("C14").Value = 2512
$ws.Range("E14").Value = 4430
$ws.Range("G14").Value = "75ms"

# Row 15 (Maze100_3.txt)
$ws.Range("C15").Value = 223
$ws.Range("E15").Value = 9821
$ws.Range("G15").Value = "72ms"

# Row 16 (Maze100_4.txt)
$ws.Range("C16").Value = 3018
$ws.Range("E16").Value = 6690
$ws.Range("G16").Value = "77ms"

# Update the active selection to match the saved view state
[void]$ws.Range("D15").Select()
